$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("FAPs","Sema3g","Nrp2","FAPs",2,0.6666666666666666,9.629280666666666,28.887842,0.8773672424304264,0.8773672424304263,3,1,34.52052333333333,103.56157,0.7684334662422598,0.7684334662422598,332.4078079368822,2991.670271431939,0.6741983512682257,0.6741983512682256),
    @("FAPs","Sema3g","Nrp2","sCs",2,0.6666666666666666,9.629280666666666,28.887842,0.8773672424304264,0.8773672424304263,3,1,5.347618333333333,16.042855,0.1190390091234806,0.1190390091234805,51.49371782987888,463.44346046891,0.1044409271763185,0.1044409271763185),
    @("FAPs","Sema3g","Nrp2","ECs",2,0.6666666666666666,9.629280666666666,28.887842,0.8773672424304264,0.8773672424304263,3,1,5.055101333333334,15.165304,0.1125275246342597,0.1125275246342597,48.67698953710756,438.092905833968,0.09872796398588231,0.09872796398588229),
    @("sCs","Sema3g","Nrp2","FAPs",3,1,0.2843916666666667,0.853175,0.02591220891683702,0.02591220891683702,3,1,34.52052333333333,103.56157,0.7684334662422598,0.7684334662422598,9.817349164972221,88.35614248474999,0.01991180851595866,0.01991180851595866),
    @("sCs","Sema3g","Nrp2","sCs",3,1,0.2843916666666667,0.853175,0.02591220891683702,0.02591220891683702,3,1,5.347618333333333,16.042855,0.1190390091234806,0.1190390091234805,1.520818090513889,13.687362814625,0.003084563673660896,0.003084563673660896),
    @("sCs","Sema3g","Nrp2","ECs",3,1,0.2843916666666667,0.853175,0.02591220891683702,0.02591220891683702,3,1,5.055101333333334,15.165304,0.1125275246342597,0.1125275246342597,1.437628693355556,12.9386582402,0.002915836727217462,0.002915836727217461),
    @("ECs","Sema3g","Nrp2","FAPs",3,1,1.061527333333333,3.184582,0.09672054865273674,0.09672054865273673,3,1,34.52052333333333,103.56157,0.7684334662422598,0.7684334662422598,36.64447907930444,329.80031171374,0.07432330645807562,0.07432330645807562),
    @("ECs","Sema3g","Nrp2","sCs",3,1,1.061527333333333,3.184582,0.09672054865273674,0.09672054865273673,3,1,5.347618333333333,16.042855,0.1190390091234806,0.1190390091234805,5.676643029067777,51.08978726160999,0.01151351827350118,0.01151351827350117),
    @("ECs","Sema3g","Nrp2","ECs",3,1,1.061527333333333,3.184582,0.09672054865273674,0.09672054865273673,3,1,5.055101333333334,15.165304,0.1125275246342597,0.1125275246342597,5.366128238103111,48.295154142928,0.01088372392115995,0.01088372392115994)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
